$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 and Row 4 swap their Fecha (date), Volumen, Precio minimo,
# Precio promedio ponderado, and Precio $/Kg values.

# New Row 3 values (previously Row 4's values)
$ws.Range("D3").Value2 = 44714
$ws.Range("J3").Value2 = 80
$ws.Range("K3").Value2 = 9000
$ws.Range("M3").Value2 = 9500
$ws.Range("P3").Value2 = 528

# New Row 4 values (previously Row 3's values)
$ws.Range("D4").Value2 = 44804
$ws.Range("J4").Value2 = 50
$ws.Range("K4").Value2 = 9500
$ws.Range("M4").Value2 = 9750
$ws.Range("P4").Value2 = 542
